# "Added last minute updates"
#
# The first paragraph of the document is an internal bookmark/ID
# paragraph ("**ID__AFFARS_5327_topic_3__ID**" + a trailing space run).
# This change:
#   1. Renames the ID text to match the actual topic
#      (5327_topic_3 -> 5327_201_2) and drops the now-unused trailing
#      space run.
#   2. Gives the paragraph a (invisible, spacing-only) paragraph border
#      with 5pt space on every side, matching the borders already used
#      on the final paragraph of the document.
#   3. Bumps the paragraph's left indent from 120 twips (6pt) to
#      225 twips (11.25pt) to match that same paragraph's indent.

$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)

# --- paragraph formatting -------------------------------------------------

# Add a paragraph border (<w:pBdr>) with w:space="5" on all four sides.
$p1.Format.Borders.DistanceFromTop = 5
$p1.Format.Borders.DistanceFromLeft = 5
$p1.Format.Borders.DistanceFromBottom = 5
$p1.Format.Borders.DistanceFromRight = 5

# <w:ind w:left="120"/> -> <w:ind w:left="225"/>  (120/20=6pt, 225/20=11.25pt)
$p1.Format.LeftIndent = 11.25

# --- run/text content ------------------------------------------------------

$r1 = $p1.Range
$paraStart = $r1.Start
$paraEnd = $r1.End

# Remove the trailing " " run entirely (last character of the paragraph,
# right before the paragraph mark), so only the id run remains.
$spaceRange = $d.Range($paraEnd - 2, $paraEnd - 1)
$spaceRange.Text = ""

# Rename the bookmark id text in the remaining run.
$idRange = $d.Range($paraStart, $paraEnd - 2)
$idRange.Text = "**ID__AFFARS_5327_201_2__ID**"
